$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0.008118172954500277
$ws.Range("E2").Value = 0.856825922154505
$ws.Range("F2").Value = 0.0009081363677978516
$ws.Range("G2").Value = 1.399548801684356

# Row 3 updates
$ws.Range("B3").Value = 4722
$ws.Range("C3").Value = 0.05103141248772665
$ws.Range("E3").Value = 0.09999756533481143
$ws.Range("F3").Value = 2.449036121368408
$ws.Range("G3").Value = 0.7429462184840422
